# Generate Report for Handback
#
# This script updates the "localization-status.xlsx" handback report:
#   - Overview sheet: "In Translation" status becomes
#     "Handed back: in sync with en-US" (and the two status columns widen
#     to fit the longer text).
#   - zh-cn / de-de sheets: the "Latest Target File" / "Latest Handback
#     File" / "Latest Handback DateTime" columns get populated now that
#     handback has happened, with a hyperlink on the target-file cell
#     (mirroring the existing source-file hyperlink).

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ed6922ac15bd2d27c23e89ea50dce599ba630e68/e2e/"
$file1 = "21796476-f596-4d61-9fb3-e89e0d3f509e.md"
$file2 = "8b6da137-1598-4b85-a629-cbc57e03b9f3.md"

# Column width helper: this host stores ColumnWidth rounded to the nearest
# 1/6 of a character, then serialises `width = ColumnWidth + 5/6` into the
# OOXML <col> element. Pick the ColumnWidth that lands on that grid.
function Set-ColWidthForXml($range, [double]$xmlWidth) {
    $sixths = [Math]::Round(($xmlWidth - 5.0/6.0) * 6.0)
    $range.ColumnWidth = $sixths / 6.0
}

# ---------------------------------------------------------------------
# Overview sheet: handback status text
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

Set-ColWidthForXml $overview.Range("E1") 29.9777047293527
Set-ColWidthForXml $overview.Range("F1") 29.9777047293527

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): fill in handback columns
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; XlfSuffix = "zh-cn"; HandbackDateTime = "2016-08-29 16:27:05" },
    @{ Sheet = "de-de"; XlfSuffix = "de-de"; HandbackDateTime = "2016-08-29 16:27:17" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    Set-ColWidthForXml $ws.Range("C1") 29.9777047293527
    Set-ColWidthForXml $ws.Range("I1") 40
    Set-ColWidthForXml $ws.Range("J1") 40

    # --- Row 2 (21796476-...md) ---
    $ws.Range("I2").Value = $file1
    $ws.Hyperlinks.Add($ws.Range("I2"), ($repoBase + $file1), "", "", $file1)
    $ws.Range("J2").Value = "21796476-f596-4d61-9fb3-e89e0d3f509e.e8cdd7c1a2c545774348e5f73e3db68a35713c2e." + $lang.XlfSuffix + ".xlf"
    $ws.Range("K2").Value = $lang.HandbackDateTime

    # --- Row 3 (8b6da137-...md) ---
    $ws.Range("I3").Value = $file2
    $ws.Hyperlinks.Add($ws.Range("I3"), ($repoBase + $file2), "", "", $file2)
    $ws.Range("J3").Value = "8b6da137-1598-4b85-a629-cbc57e03b9f3.6bbd60c6035fa9f8510da076da74b72fd27f60e1." + $lang.XlfSuffix + ".xlf"
    $ws.Range("K3").Value = $lang.HandbackDateTime
}
